$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column R (18) -- this pushes the old
# (empty) R column to T, and the old S "Resources" column to U, matching
# the target layout where two new week columns are added after Wk16 (Q).
$ws.Range("R1:S1").EntireColumn.Insert()

# New week header cells (Wk 17, Wk 18, Wk 19), matching the style used by
# the other week headers J3:Q3.
$ws.Range("O3:Q3").Copy()
$ws.Range("R3:T3").PasteSpecial(-4122)
$ws.Range("R3").Value = "Wk 17 10/09"
$ws.Range("S3").Value = "Wk 18 17/09"
$ws.Range("T3").Value = "Wk 19 24/09"

# "Completed" markers for Wk16 (existing column Q) on rows 12-14.
$ws.Range("Q12").Value = "Completed"
$ws.Range("Q13").Value = "Completed"
$ws.Range("Q14").Value = "Completed"
$ws.Range("Q12:Q14").Font.Bold = $true
$ws.Range("Q12:Q14").HorizontalAlignment = -4108

# "Completed" markers for the new Wk18 column (S) on rows 16-19.
$ws.Range("S16").Value = "Completed"
$ws.Range("S17").Value = "Completed"
$ws.Range("S18").Value = "Completed"
$ws.Range("S19").Value = "Completed"
$ws.Range("S16:S19").Font.Bold = $true

# Column widths: Q (17) and S (19) get an explicit width of 10.
$ws.Columns(17).ColumnWidth = 9.14
$ws.Columns(19).ColumnWidth = 9.14

# Minor row-height tweaks on rows 12-14 (as recorded by Excel for this edit).
$ws.Rows(12).RowHeight = 13.8
$ws.Rows(13).RowHeight = 15.6
$ws.Rows(14).RowHeight = 15

# Selection, to match the saved cursor position.
$ws.Range("S20").Select()
